$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999440806064,
    0.99465453603624021,
    0.97351049367626885,
    0.96420467410530153,
    0.95532684232830645,
    0.93627242267728672,
    0.93516880254849744,
    0.92937748102826845,
    0.9257427499046349,
    0.92303316097318078,
    0.92223595444528184,
    0.92097507308957938,
    0.91598214917179299,
    0.91427526534386727,
    0.91422721504638005,
    0.91172075774321915,
    0.90801286395284841,
    0.90690395059458984,
    0.99249922090809406,
    0.9752676968165922,
    0.9713885111600693,
    0.97012400096177032,
    0.96599823133845208,
    0.95297651215265211,
    0.94651937564088628,
    0.9341701238146749,
    0.93148854236292278,
    0.9209668011581893,
    0.91382381906539389,
    0.91097013822392725,
    0.91388986473403122,
    0.91598789474738362,
    0.92060925825521311
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
